$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in B, C, D columns for rows 281-302 ---
$ws.Range("B281").Value = 277
$ws.Range("C281").Value = 43531
$ws.Range("D281").Value = 'Likeshuo'
$ws.Range("B282").Value = 278
$ws.Range("C282").Value = 43531
$ws.Range("D282").Value = 'Likeshuo'
$ws.Range("B283").Value = 279
$ws.Range("C283").Value = 43531
$ws.Range("D283").Value = 'Likeshuo'
$ws.Range("B284").Value = 280
$ws.Range("C284").Value = 43531
$ws.Range("D284").Value = 'Likeshuo'
$ws.Range("B285").Value = 281
$ws.Range("C285").Value = 43531
$ws.Range("D285").Value = 'Likeshuo'
$ws.Range("B286").Value = 282
$ws.Range("C286").Value = 43531
$ws.Range("D286").Value = 'Likeshuo'
$ws.Range("B287").Value = 283
$ws.Range("C287").Value = 43531
$ws.Range("D287").Value = 'Likeshuo'
$ws.Range("B288").Value = 284
$ws.Range("C288").Value = 43531
$ws.Range("D288").Value = 'Likeshuo'
$ws.Range("B289").Value = 285
$ws.Range("C289").Value = 43531
$ws.Range("D289").Value = 'Likeshuo'
$ws.Range("B290").Value = 286
$ws.Range("C290").Value = 43531
$ws.Range("D290").Value = 'Likeshuo'
$ws.Range("B291").Value = 287
$ws.Range("C291").Value = 43531
$ws.Range("D291").Value = 'Likeshuo'
$ws.Range("B292").Value = 288
$ws.Range("C292").Value = 43531
$ws.Range("D292").Value = 'Likeshuo'
$ws.Range("B293").Value = 289
$ws.Range("C293").Value = 43531
$ws.Range("D293").Value = 'Likeshuo'
$ws.Range("B294").Value = 290
$ws.Range("C294").Value = 43531
$ws.Range("D294").Value = 'Likeshuo'
$ws.Range("B295").Value = 291
$ws.Range("C295").Value = 43531
$ws.Range("D295").Value = 'Likeshuo'
$ws.Range("B296").Value = 292
$ws.Range("C296").Value = 43531
$ws.Range("D296").Value = 'Likeshuo'
$ws.Range("B297").Value = 293
$ws.Range("C297").Value = 43531
$ws.Range("D297").Value = 'Likeshuo'
$ws.Range("B298").Value = 294
$ws.Range("C298").Value = 43531
$ws.Range("D298").Value = 'Likeshuo'
$ws.Range("B299").Value = 295
$ws.Range("C299").Value = 43531
$ws.Range("D299").Value = 'Likeshuo'
$ws.Range("B300").Value = 296
$ws.Range("C300").Value = 43531
$ws.Range("D300").Value = 'Likeshuo'
$ws.Range("B301").Value = 297
$ws.Range("C301").Value = 43531
$ws.Range("D301").Value = 'Likeshuo'
$ws.Range("B302").Value = 298
$ws.Range("C302").Value = 43531
$ws.Range("D302").Value = 'Likeshuo'

# --- Fill in E column (English words) in row order, so shared strings are appended in this order ---
$ws.Range("E281").Value = 'fraud'
$ws.Range("E282").Value = 'coercion'
$ws.Range("E283").Value = 'deceive'
$ws.Range("E284").Value = 'burglar'
$ws.Range("E285").Value = 'burglary'
$ws.Range("E286").Value = 'intent'
$ws.Range("E287").Value = 'plainclothes police officer'
$ws.Range("E288").Value = 'charge'
$ws.Range("E289").Value = 'impetuous'
$ws.Range("E290").Value = 'offense'
$ws.Range("E291").Value = 'stab'
$ws.Range("E292").Value = 'stabbed'
$ws.Range("E293").Value = 'strict'
$ws.Range("E294").Value = 'rigidly'
$ws.Range("E295").Value = 'deviation'
$ws.Range("E296").Value = 'follow the tracks of'
$ws.Range("E297").Value = 'fine'
$ws.Range("E298").Value = 'something you have to pay'
$ws.Range("E299").Value = 'smuggler'
$ws.Range("E300").Value = 'lucrative'
$ws.Range("E301").Value = 'covert'
$ws.Range("E302").Value = 'hidden'

# --- Fill in G column (Chinese translations) in row order ---
$ws.Range("G281").Value = '骗局'
$ws.Range("G282").Value = '强迫'
$ws.Range("G283").Value = '欺骗'
$ws.Range("G284").Value = '窃贼'
$ws.Range("G285").Value = '窃案'
$ws.Range("G286").Value = '意图'
$ws.Range("G287").Value = '便衣警察'
$ws.Range("G288").Value = '收费'
$ws.Range("G289").Value = '浮躁'
$ws.Range("G290").Value = '罪行'
$ws.Range("G291").Value = '刺'
$ws.Range("G292").Value = '被刺'
$ws.Range("G293").Value = '严格'
$ws.Range("G294").Value = '严格'
$ws.Range("G295").Value = '偏差'
$ws.Range("G296").Value = '按照轨道'
$ws.Range("G297").Value = '精细'
$ws.Range("G298").Value = '你需要支付的东西'
$ws.Range("G299").Value = '走私者'
$ws.Range("G300").Value = '有利可图'
$ws.Range("G301").Value = '隐蔽'
$ws.Range("G302").Value = '隐'

# --- Fill in K and L columns ---
$ws.Range("K281").Value = 1
$ws.Range("L281").Value = 'Input'
$ws.Range("K282").Value = 1
$ws.Range("L282").Value = 'Input'
$ws.Range("K283").Value = 1
$ws.Range("L283").Value = 'Input'
$ws.Range("K284").Value = 1
$ws.Range("L284").Value = 'Input'
$ws.Range("K285").Value = 1
$ws.Range("L285").Value = 'Input'
$ws.Range("K286").Value = 1
$ws.Range("L286").Value = 'Input'
$ws.Range("K287").Value = 1
$ws.Range("L287").Value = 'Input'
$ws.Range("K288").Value = 1
$ws.Range("L288").Value = 'Input'
$ws.Range("K289").Value = 1
$ws.Range("L289").Value = 'Input'
$ws.Range("K290").Value = 1
$ws.Range("L290").Value = 'Input'
$ws.Range("K291").Value = 1
$ws.Range("L291").Value = 'Input'
$ws.Range("K292").Value = 1
$ws.Range("L292").Value = 'Input'
$ws.Range("K293").Value = 1
$ws.Range("L293").Value = 'Input'
$ws.Range("K294").Value = 1
$ws.Range("L294").Value = 'Input'
$ws.Range("K295").Value = 1
$ws.Range("L295").Value = 'Input'
$ws.Range("K296").Value = 1
$ws.Range("L296").Value = 'Input'
$ws.Range("K297").Value = 1
$ws.Range("L297").Value = 'Input'
$ws.Range("K298").Value = 1
$ws.Range("L298").Value = 'Input'
$ws.Range("K299").Value = 1
$ws.Range("L299").Value = 'Input'
$ws.Range("K300").Value = 1
$ws.Range("L300").Value = 'Input'
$ws.Range("K301").Value = 1
$ws.Range("L301").Value = 'Input'
$ws.Range("K302").Value = 1
$ws.Range("L302").Value = 'Input'

# --- Update data validation ranges to reflect the new rows (281-302) ---
$rOldA = $ws.Range("D281:D1048576")
$rOldA.Validation.Delete()
$rNewA = $ws.Range("D303:D1048576")
$rNewA.Validation.Add(3, 1, 1, '"Likeshuo,TOEFL,TPO"')

$rOldB = $ws.Range("D258:D280")
$rOldB.Validation.Delete()
$rNewB = $ws.Range("D258:D302")
$rNewB.Validation.Add(3, 1, 1, '"Likeshuo,TOEFL,TPO, 500 setns, NCE4"')

# --- Update the active view/selection to match the final cursor position ---
$ws.Activate()
$ws.Range("F294").Select()
